$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Drop the now-unused trailing columns U:AD (row1 used 0..28, now only 0..18 needed)
$ws.Range("U1:AD2").EntireColumn.Delete()

# 2. Row 2 header labels (C2:T2) - Miller index / pair labels
$ws.Range("C2").Value = "[3, 2, 1]"
$ws.Range("D2").Value = "[1, 1, 0]"
$ws.Range("E2").Value = "[2, 2, 2]"
$ws.Range("F2").Value = "[3, 1, 0]"
$ws.Range("G2").Value = "[2, 2, 0]"
$ws.Range("H2").Value = "[2, 0, 0]"
$ws.Range("I2").Value = "[2, 1, 1]"
$ws.Range("J2").Value = "[4, 0, 0]"
$ws.Range("K2").Value = "1Pair-A"
$ws.Range("L2").Value = "1Pair-B"
$ws.Range("M2").Value = "2Pairs-A"
$ws.Range("N2").Value = "2Pairs-B"
$ws.Range("O2").Value = "3Pairs-A"
$ws.Range("P2").Value = "3Pairs-B"
$ws.Range("Q2").Value = "3Pairs-C"
$ws.Range("R2").Value = "4Pairs"
$ws.Range("S2").Value = "5A4F"
$ws.Range("T2").Value = "MaxUnique"

# 3. Relabel the scheme rows in column B (rows 3-19 keep their row position, new text)
$schemeNames = @(
  "Spiral5",
  "RotRing OmegaMax-90",
  "Equal Angle",
  "Tilt Rotate",
  "CLR",
  "Rizzie Hex",
  "Thomas Hex",
  "Tilt Rotate_Partial",
  "RotRing OmegaMax-60",
  "Equal Angle_Partial",
  "Rizzie Hex_Partial",
  "ND Single",
  "RD Single",
  "TD Single",
  "Morris Single",
  "Ring Perpendicular to ND",
  "Ring Perpendicular to RD",
  "Ring Perpendicular to TD",
  "OffsetFTD",
  "OffsetATD",
  "OffsetF45",
  "OffsetA45",
  "OffsetFRD",
  "OffsetARD",
  "Gaussian Quadrature",
  "Michael-CCHex",
  "Michael-SNHex"
)

# 4. New rows 20-29 need to carry the same "A column" formatting (bold/border/
#    center) as the pre-existing rows, so clone it down from row 19 first.
$ws.Range("A19").Copy() | Out-Null
$ws.Range("A20:A29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $schemeNames.Count; $i++) {
    $r = 3 + $i
    $ws.Cells.Item($r, 1).Value = $i + 1
    $ws.Cells.Item($r, 2).Value = $schemeNames[$i]
    for ($c = 3; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}
